# Updated Excel examples in conversion
#
# The "Events" and "Locations" example sheets referenced the vocabulary
# namespaces using a hash ("#") separator; update them to use a trailing
# slash ("/") separator instead, and leave the workbook with the
# "Locations" sheet/tab active (matching the last-saved selection state).

$wb = $excel.ActiveWorkbook

$wsEvents    = $wb.Worksheets.Item("Events")
$wsLocations = $wb.Worksheets.Item("Locations")

# --- Events sheet: update the two namespace URIs from '#' to '/' ---
$wsEvents.Range("C2").Value = "http://data.sparna.fr/vocabularies/semweb-events/"
$wsEvents.Range("C3").Value = "http://data.sparna.fr/vocabularies/places/"

# --- Locations sheet: update the namespace URI from '#' to '/' ---
$wsLocations.Range("C2").Value = "http://data.sparna.fr/vocabularies/places/"

# --- Update selections / active sheet to match the saved workbook state ---
$wsEvents.Range("C4").Select()

$wsLocations.Activate()
$wsLocations.Range("C3").Select()
